# XLS and XLSX conversions now use agate-excel.
#
# Before: two sheets -- "testfixed_converted.csv" (the real data, incl. a
# "time" column) and "Sheet2" (a near-duplicate of the data, tab-selected).
# After: the duplicate sheet is wiped and renamed to "not this one", and the
# real data sheet loses its redundant "time" column (F), gets renamed to
# "data", and becomes the selected/active sheet with C4 selected.

$wb = $excel.ActiveWorkbook

$wsData  = $wb.Worksheets.Item("testfixed_converted.csv")
$wsExtra = $wb.Worksheets.Item("Sheet2")

# The "time" column duplicates the "datetime" column as a day fraction --
# drop it, shifting "datetime"/"empty_column" left.
$null = $wsData.Columns.Item(6).Delete()
$null = $wsData.Unprotect()
$wsData.PageSetup.LeftMargin = 54
$wsData.PageSetup.RightMargin = 54
$wsData.PageSetup.TopMargin = 72
$wsData.PageSetup.BottomMargin = 72
$wsData.PageSetup.HeaderMargin = 36
$wsData.PageSetup.FooterMargin = 36

# The second sheet was just a re-ordered copy of the same data -- clear it
# out entirely, it becomes a blank placeholder sheet.
$null = $wsExtra.Cells.Clear()
$null = $wsExtra.Unprotect()
$wsExtra.PageSetup.LeftMargin = 54
$wsExtra.PageSetup.RightMargin = 54
$wsExtra.PageSetup.TopMargin = 72
$wsExtra.PageSetup.BottomMargin = 72
$wsExtra.PageSetup.HeaderMargin = 36
$wsExtra.PageSetup.FooterMargin = 36

# Rename both sheets.
$wsExtra.Name = "not this one"
$wsData.Name = "data"

# Put the (now blank) "not this one" sheet first, "data" second -- matches
# the original rId1/rId2 physical ordering.
$null = $wsExtra.Move($wb.Worksheets.Item(1))

# Re-fetch by name: Move() invalidates old positional handles.
$wsData = $wb.Worksheets.Item("data")
$null = $wsData.Activate()
$null = $wsData.Range("C4").Select()
